# New crime data collected — refresh the weekly CompStat report:
#   * bump the report Volume/Number and the "week covering" date range
#   * replace the week's crime-complaint figures (rows 15-27) with the
#     newly collected counts/percentages

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Masthead: Volume 30 Number 45 -> 46; week of 11/6-11/12 -> 11/13-11/19
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# ---------------------------------------------------------------------
# A few cells flip between the "***.*"/"0" placeholder text and a real
# number (or vice versa) this week, so they need their number format
# set explicitly before the new numeric value is written — otherwise
# the cell keeps rendering as the old placeholder text style. For the
# lone cell that goes the other way (number -> text placeholder "0"),
# copy the format+value from another cell that already shows that
# exact placeholder.
# ---------------------------------------------------------------------
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1

$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1

$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1

$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = 0

$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 1

$ws.Range("D15").Copy($ws.Range("C27"))

# ---------------------------------------------------------------------
# Row 15 — Rape
# ---------------------------------------------------------------------
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = -37.5
$ws.Range("L15").Value = 11.111111111111
$ws.Range("M15").Value = 11.111111111111
$ws.Range("N15").Value = -44.444444444444

# ---------------------------------------------------------------------
# Row 16 — Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -22.222222222222
$ws.Range("I16").Value = 121
$ws.Range("J16").Value = 166
$ws.Range("K16").Value = -27.108433734939
$ws.Range("L16").Value = 22.222222222222
$ws.Range("M16").Value = -3.2
$ws.Range("N16").Value = -83.401920438957

# ---------------------------------------------------------------------
# Row 17 — Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 194
$ws.Range("J17").Value = 194
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8.379888268156
$ws.Range("M17").Value = 61.666666666666
$ws.Range("N17").Value = -5.825242718446

# ---------------------------------------------------------------------
# Row 18 — Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 155
$ws.Range("J18").Value = 147
$ws.Range("K18").Value = 5.442176870748
$ws.Range("L18").Value = 55
$ws.Range("M18").Value = 76.136363636363
$ws.Range("N18").Value = -55.202312138728

# ---------------------------------------------------------------------
# Row 19 — Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 533
$ws.Range("J19").Value = 693
$ws.Range("K19").Value = -23.088023088023
$ws.Range("L19").Value = 8.113590263691
$ws.Range("M19").Value = 111.507936507937
$ws.Range("N19").Value = 30.317848410757

# ---------------------------------------------------------------------
# Row 20 — G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 59
$ws.Range("K20").Value = 5.357142857142
$ws.Range("L20").Value = 13.461538461538
$ws.Range("M20").Value = 28.260869565217
$ws.Range("N20").Value = -82.947976878612

# ---------------------------------------------------------------------
# Row 21 — TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -3.333333333333
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = 13.186813186813
$ws.Range("I21").Value = 1074
$ws.Range("J21").Value = 1275
$ws.Range("K21").Value = -15.764705882352
$ws.Range("L21").Value = 14.621131270010
$ws.Range("M21").Value = 67.550702028081
$ws.Range("N21").Value = -47.914645974781

# ---------------------------------------------------------------------
# Row 22 — Transit
# ---------------------------------------------------------------------
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 17
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -29.166666666666
$ws.Range("L22").Value = 21.428571428571
$ws.Range("M22").Value = 6.25

# ---------------------------------------------------------------------
# Row 23 — Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 61.538461538461
$ws.Range("I23").Value = 155
$ws.Range("J23").Value = 160
$ws.Range("K23").Value = -3.125
$ws.Range("L23").Value = 4.026845637583
$ws.Range("M23").Value = 27.049180327868

# ---------------------------------------------------------------------
# Row 24 — Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 4.166666666666
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = -13.934426229508
$ws.Range("I24").Value = 1115
$ws.Range("J24").Value = 2020
$ws.Range("K24").Value = -44.801980198019
$ws.Range("L24").Value = -20.071684587813
$ws.Range("M24").Value = 65.676077265973

# ---------------------------------------------------------------------
# Row 25 — Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 20
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 2.702702702702
$ws.Range("I25").Value = 404
$ws.Range("J25").Value = 404
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 39.792387543252
$ws.Range("M25").Value = 40.277777777777

# ---------------------------------------------------------------------
# Row 26 — UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 19
$ws.Range("K26").Value = -24
$ws.Range("L26").Value = 35.714285714285

# ---------------------------------------------------------------------
# Row 27 — Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = 4.761904761904
$ws.Range("L27").Value = 7.317073170731
